$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Rebuild the "Key Terms (Biblica) is based on: ..." paragraph in place.
#    It becomes the new "Biblica Study Notes (Key Terms) ..." resource blurb.
#    We find it by its distinctive bold lead-in run text.
# ---------------------------------------------------------------------------
$targetPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t.StartsWith("Key Terms (Biblica) is based on")) {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    $r = $targetPara.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = ""

    # Bold run: "Biblica Study Notes (Key Terms)"
    $r1 = $targetPara.Range
    $r1.MoveEnd(1, -1) | Out-Null
    $r1.Collapse(0) | Out-Null
    $r1.Font.Bold = $true
    $r1.InsertAfter("Biblica Study Notes (Key Terms)")

    # Normal run: " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
    $r2 = $targetPara.Range
    $r2.MoveEnd(1, -1) | Out-Null
    $r2.Collapse(0) | Out-Null
    $r2.Font.Bold = $false
    $r2.InsertAfter(" © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. ")

    # Normal run: "Biblica Study Notes"
    $r3 = $targetPara.Range
    $r3.MoveEnd(1, -1) | Out-Null
    $r3.Collapse(0) | Out-Null
    $r3.Font.Bold = $false
    $r3.InsertAfter("Biblica Study Notes")

    # Normal run: " has been adapted in the following languages: ..."
    $r4 = $targetPara.Range
    $r4.MoveEnd(1, -1) | Out-Null
    $r4.Collapse(0) | Out-Null
    $r4.Font.Bold = $false
    $r4.InsertAfter(" has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual.")
}

# ---------------------------------------------------------------------------
# 2) Delete the now-obsolete paragraphs entirely (including their paragraph
#    marks): "License Information" heading, "This PDF version is provided
#    under the same license." and the "Gabriel, Gad, Galatia, ..." summary
#    list paragraph. Collect matches first, then delete from the bottom of
#    the document upward so earlier indices stay valid.
# ---------------------------------------------------------------------------
$toDelete = New-Object System.Collections.ArrayList

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "License Information") {
        $toDelete.Add($i) | Out-Null
    } elseif ($t -eq "This PDF version is provided under the same license.") {
        $toDelete.Add($i) | Out-Null
    } elseif ($t.StartsWith("Gabriel, Gad, Galatia")) {
        $toDelete.Add($i) | Out-Null
    }
}

$sorted = $toDelete | Sort-Object -Descending
foreach ($idx in $sorted) {
    $p = $d.Paragraphs($idx)
    $p.Range.Delete()
}
